# Apply the re-ordering of rows 2-21 (columns A:F) on the active sheet.
# The underlying data set is unchanged; only the row order differs
# (this mirrors a re-shuffle of the componentsMapping / symbol-weight
# table used by lowcode.Config).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1203, 3,  15, 15, 15, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(901,  16, 15, 45, 60, 60),
    @(902,  1,  0,  0,  0,  0),
    @(401,  9,  48, 67, 75, 45),
    @(801,  3,  67, 65, 52, 45),
    @(1201, 2,  10, 10, 10, 10),
    @(1202, 2,  10, 10, 10, 10),
    @(101,  9,  30, 15, 60, 15),
    @(501,  9,  52, 30, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(601,  9,  60, 67, 60, 42),
    @(201,  9,  30, 15, 45, 30),
    @(301,  6,  45, 30, 60, 45),
    @(1101, 0,  15, 30, 30, 0),
    @(802,  0,  4,  5,  4,  0),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(1,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
